# Week 13 logging update: update target depth data on row "H" (home) for
# both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 442
$wsOff.Range("C2").Value = 321
$wsOff.Range("D2").Value = 112
$wsOff.Range("E2").Value = 47

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 441
$wsDef.Range("C2").Value = 292
